# Fruta / hortaliza, semanal
# Insert two new weekly price rows (2022-12-23, date serial 44918) for
# "Feria Lagunitas de Puerto Montt" / Frutilla, right after the existing
# row 300, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 301:302 - everything that was in 301:329
# shifts down to 303:331.
$ws.Rows("301:302").Insert()

# --- New row 301 --------------------------------------------------------
$ws.Range("A301").Value = 4
$ws.Range("B301").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C301").Value = "Los Lagos"
$ws.Range("D301").Value = 44918
$ws.Range("E301").Value = 10
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100101
$ws.Range("H301").Value = "Berries"
$ws.Range("I301").Value = 100112025
$ws.Range("J301").Value = "Frutilla"
$ws.Range("K301").Value = "Sin especificar"
$ws.Range("L301").Value = "Primera"
$ws.Range("M301").Value = 800
$ws.Range("N301").Value = 10000
$ws.Range("O301").Value = 11000
$ws.Range("P301").Value = 10500
$ws.Range("Q301").Value = "$/bandeja 7 kilos"
$ws.Range("R301").Value = "Provincia de Melipilla"
$ws.Range("S301").Value = 1500
$ws.Range("T301").Value = 7

# --- New row 302 --------------------------------------------------------
$ws.Range("A302").Value = 4
$ws.Range("B302").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C302").Value = "Los Lagos"
$ws.Range("D302").Value = 44918
$ws.Range("E302").Value = 10
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100101
$ws.Range("H302").Value = "Berries"
$ws.Range("I302").Value = 100112025
$ws.Range("J302").Value = "Frutilla"
$ws.Range("K302").Value = "Sin especificar"
$ws.Range("L302").Value = "Primera"
$ws.Range("M302").Value = 800
$ws.Range("N302").Value = 11000
$ws.Range("O302").Value = 12000
$ws.Range("P302").Value = 11500
$ws.Range("Q302").Value = "$/caja 7 kilos"
$ws.Range("R302").Value = "Región de La Araucanía"
$ws.Range("S302").Value = 1643
$ws.Range("T302").Value = 7
